$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M2").Value = 784.0599999999999
$ws1.Range("M16").Value = -7825.69

# --- Sheet 2: "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F2").Value = 784.0599999999999
$ws2.Range("F16").Value = -7825.69
$ws2.Range("F38").Value = 1714.32

# --- Sheet 3: "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D12").Value = 1028.07
$ws3.Range("E12").Value = 30185.93
$ws3.Range("F12").Value = 0.03293618248221952

$ws3.Range("D14").Value = 2473.98
$ws3.Range("E14").Value = 37805.58164865473
$ws3.Range("F14").Value = 0.06142023147073217

# Widen column F on "CUMPLIMIENTO MENSUAL" from 24 to 25 (stored OOXML
# character-width units). Excel's `ColumnWidth` COM property is expressed
# in "Normal style" character widths that get re-quantized to the internal
# Maximum-Digit-Width pixel grid on write (stored = (round(ColumnWidth*6)+5)/6),
# so 24.17 (which lands in the same pixel bucket as a raw width of 25) is
# used here to land exactly on width="25" once persisted.
$ws3.Columns("F").ColumnWidth = 24.17
